# "for Code Review 3" - fill in the Code Review 3 (column D) and Code Review 4 (column E)
# marks for each student, matching the bottom-border / header-border styling used
# elsewhere on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel COM constants used below:
#   xlContinuous = 1          (LineStyle: solid line)
#   xlNone       = -4142      (LineStyle: no line)
#   xlThin       = 2          (BorderWeight)
#   xlEdgeLeft   = 7
#   xlEdgeTop    = 8
#   xlEdgeBottom = 9
#   xlEdgeRight  = 10

# ---------------------------------------------------------------------------
# 1. Fill in the Code Review 3 (D) and Code Review 4 (E) marks
# ---------------------------------------------------------------------------
$ws.Range("D3").Value = 20
$ws.Range("E3").Value = 20

$ws.Range("D4").Value = 20
$ws.Range("E4").Value = 20

$ws.Range("D5").Value = 20
$ws.Range("E5").Value = 20

$ws.Range("D6").Value = 20
$ws.Range("E6").Value = 20

$ws.Range("D7").Value = 20
$ws.Range("E7").Value = 20

# ---------------------------------------------------------------------------
# 2. Fix up the borders so columns D and E pick up the same "table" look that
#    columns B/C/F/G already have: a separator line under the header row and
#    another separator line above the totals row.
# ---------------------------------------------------------------------------

# Header cells D2:E2 lose the border under them (it now continues down into
# the data rows), keeping their grey fill and the thin box around the top.
$headerDE = $ws.Range("D2:E2")
$headerDE.Borders.Item(9).LineStyle = -4142

# A thin line above the totals row, which also works as the line below the
# last data row (row 7).
$aboveTotals = $ws.Range("D4:E7")
$aboveTotals.Borders.Item(9).LineStyle = 1
$aboveTotals.Borders.Item(9).Weight = 2

# The totals row D8:E8 no longer needs a line on top (row 7 already supplies
# it), just keep the line underneath.
$totalsDE = $ws.Range("D8:E8")
$totalsDE.Borders.Item(8).LineStyle = -4142

# ---------------------------------------------------------------------------
# 3. Misc. view state that Excel records when a user re-saves the workbook.
# ---------------------------------------------------------------------------
$ws.Range("K12").Select()
